$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025617356788274
$ws.Range("D2").Value = 1.028261897521059
$ws.Range("E2").Value = 1.025942288655465
$ws.Range("F2").Value = 1.035578244822902
$ws.Range("I2").Value = 1.031350482319003
$ws.Range("J2").Value = 1.030785377314057
$ws.Range("K2").Value = 1.031079227073046
$ws.Range("L2").Value = 1.028766384590913
$ws.Range("M2").Value = 1.038374448093185
$ws.Range("N2").Value = 1.032249210782583

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026638059890754
$ws.Range("D3").Value = 1.029195495130404
$ws.Range("E3").Value = 1.026810832799274
$ws.Range("F3").Value = 1.036793044735852
$ws.Range("I3").Value = 1.031508602293804
$ws.Range("J3").Value = 1.031445123375345
$ws.Range("K3").Value = 1.031820383910679
$ws.Range("L3").Value = 1.029442172498026
$ws.Range("M3").Value = 1.039397594799872
$ws.Range("N3").Value = 1.032909893758947

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027298645951825
$ws.Range("D4").Value = 1.029800006832352
$ws.Range("E4").Value = 1.027373315582281
$ws.Range("F4").Value = 1.037579330063497
$ws.Range("I4").Value = 1.03160943914771
$ws.Range("J4").Value = 1.031871584866223
$ws.Range("K4").Value = 1.032299752458712
$ws.Range("L4").Value = 1.029879293839758
$ws.Range("M4").Value = 1.040059307347676
$ws.Range("N4").Value = 1.033336960874057

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027576385785985
$ws.Range("D5").Value = 1.030054241626095
$ws.Range("E5").Value = 1.027609896860112
$ws.Range("F5").Value = 1.037909939863129
$ws.Range("I5").Value = 1.031651477048833
$ws.Range("J5").Value = 1.032050764047642
$ws.Range("K5").Value = 1.032501228495845
$ws.Range("L5").Value = 1.030063021285689
$ws.Range("M5").Value = 1.040337412368454
$ws.Range("N5").Value = 1.033516394510467

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027623021254428
$ws.Range("D6").Value = 1.030096934504152
$ws.Range("E6").Value = 1.027649626537018
$ws.Range("F6").Value = 1.037965453976081
$ws.Range("I6").Value = 1.03165851463704
$ws.Range("J6").Value = 1.03208084284141
$ws.Range("K6").Value = 1.032535054237493
$ws.Range("L6").Value = 1.030093867671102
$ws.Range("M6").Value = 1.040384102810893
$ws.Range("N6").Value = 1.033546516019572

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027302357009302
$ws.Range("D7").Value = 1.029803403545333
$ws.Range("E7").Value = 1.027376476347602
$ws.Range("F7").Value = 1.037583747469032
$ws.Range("I7").Value = 1.031610002251127
$ws.Range("J7").Value = 1.031873979481064
$ws.Range("K7").Value = 1.032302444790254
$ws.Range("L7").Value = 1.02988174896613
$ws.Range("M7").Value = 1.040063023710115
$ws.Range("N7").Value = 1.033339358889527

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025962282751718
$ws.Range("D8").Value = 1.02857732569412
$ws.Range("E8").Value = 1.026235718166729
$ws.Range("F8").Value = 1.035988745257806
$ws.Range("I8").Value = 1.031404225458264
$ws.Range("J8").Value = 1.031008432167049
$ws.Range("K8").Value = 1.031329747662286
$ws.Range("L8").Value = 1.028994802887146
$ws.Range("M8").Value = 1.03872029431021
$ws.Range("N8").Value = 1.032472582399046

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02360184008755
$ws.Range("D9").Value = 1.02641999266132
$ws.Range("E9").Value = 1.024229232861523
$ws.Range("F9").Value = 1.033179864880948
$ws.Range("I9").Value = 1.031030315915033
$ws.Range("J9").Value = 1.029479882712239
$ws.Range("K9").Value = 1.029614141657628
$ws.Range("L9").Value = 1.027430692644398
$ws.Range("M9").Value = 1.036351669830141
$ws.Range("N9").Value = 1.030941862228684

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022028837579743
$ws.Range("D10").Value = 1.024983928846826
$ws.Range("E10").Value = 1.022894080607884
$ws.Range("F10").Value = 1.031308383342866
$ws.Range("I10").Value = 1.030773455160596
$ws.Range("J10").Value = 1.028458616018615
$ws.Range("K10").Value = 1.028469353853761
$ws.Range("L10").Value = 1.026387168041775
$ws.Range("M10").Value = 1.034770828686702
$ws.Range("N10").Value = 1.029919145219213

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021347856375988
$ws.Range("D11").Value = 1.024362613778711
$ws.Range("E11").Value = 1.022316545287498
$ws.Range("F11").Value = 1.030498261478656
$ws.Range("I11").Value = 1.030660434315639
$ws.Range("J11").Value = 1.028015868947512
$ws.Range("K11").Value = 1.027973402078396
$ws.Range("L11").Value = 1.025935128131649
$ws.Range("M11").Value = 1.034085881446444
$ws.Range("N11").Value = 1.029475769396486

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021094930256784
$ws.Range("D12").Value = 1.024131906604081
$ws.Range("E12").Value = 1.022102112693703
$ws.Range("F12").Value = 1.030197381369558
$ws.Range("I12").Value = 1.030618183287957
$ws.Range("J12").Value = 1.02785133318434
$ws.Range("K12").Value = 1.027789145949684
$ws.Range("L12").Value = 1.025767192540968
$ws.Range("M12").Value = 1.033831395702772
$ws.Range("N12").Value = 1.029310999973657

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021149182839643
$ws.Range("D13").Value = 1.0241813905916
$ws.Range("E13").Value = 1.022148105145147
$ws.Range("F13").Value = 1.030261919598663
$ws.Range("I13").Value = 1.030627258495032
$ws.Range("J13").Value = 1.027886630272755
$ws.Range("K13").Value = 1.027828671211614
$ws.Range("L13").Value = 1.025803216551356
$ws.Range("M13").Value = 1.033885986751751
$ws.Range("N13").Value = 1.029346347187987

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021326949005159
$ws.Range("D14").Value = 1.024343541881229
$ws.Range("E14").Value = 1.022298818389246
$ws.Range("F14").Value = 1.030473389923017
$ws.Range("I14").Value = 1.030656947340699
$ws.Range("J14").Value = 1.028002269996968
$ws.Range("K14").Value = 1.0279581721745
$ws.Range("L14").Value = 1.02592124708993
$ws.Range("M14").Value = 1.034064846923535
$ws.Range("N14").Value = 1.029462151133872

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021436479330184
$ws.Range("D15").Value = 1.02444345884149
$ws.Range("E15").Value = 1.022391689697356
$ws.Range("F15").Value = 1.030603688402588
$ws.Range("I15").Value = 1.030675203837901
$ws.Range("J15").Value = 1.02807350887851
$ws.Range("K15").Value = 1.028037957020144
$ws.Range("L15").Value = 1.025993965918036
$ws.Range("M15").Value = 1.034175039828501
$ws.Range("N15").Value = 1.029533491182797

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022074034967758
$ws.Range("D16").Value = 1.025025174251359
$ws.Range("E16").Value = 1.022932422285842
$ws.Range("F16").Value = 1.031362153464543
$ws.Range("I16").Value = 1.030780918102851
$ws.Range("J16").Value = 1.028487988488063
$ws.Range("K16").Value = 1.02850226326366
$ws.Range("L16").Value = 1.026417164524722
$ws.Range("M16").Value = 1.034816277192482
$ws.Range("N16").Value = 1.029948559400937

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022473993846238
$ws.Range("D17").Value = 1.025390205708292
$ws.Range("E17").Value = 1.023271769002334
$ws.Range("F17").Value = 1.031837982546786
$ws.Range("I17").Value = 1.030846748474416
$ws.Range("J17").Value = 1.028747838287712
$ws.Range("K17").Value = 1.028793443195486
$ws.Range("L17").Value = 1.026682575965209
$ws.Range("M17").Value = 1.035218392032909
$ws.Range("N17").Value = 1.030208778217107

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022707296550281
$ws.Range("D18").Value = 1.025603171574205
$ws.Range("E18").Value = 1.023469761524286
$ws.Range("F18").Value = 1.032115549024811
$ws.Range("I18").Value = 1.030884972683164
$ws.Range("J18").Value = 1.028899352955581
$ws.Range("K18").Value = 1.028963259284452
$ws.Range("L18").Value = 1.026837367922071
$ws.Range("M18").Value = 1.035452897134029
$ws.Range("N18").Value = 1.030360508053183

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022786849035022
$ws.Range("D19").Value = 1.025675795766444
$ws.Range("E19").Value = 1.02353728156749
$ws.Range("F19").Value = 1.032210196028625
$ws.Range("I19").Value = 1.030897976711973
$ws.Range("J19").Value = 1.028951006832817
$ws.Range("K19").Value = 1.029021158065878
$ws.Range("L19").Value = 1.02689014491112
$ws.Range("M19").Value = 1.035532850269199
$ws.Range("N19").Value = 1.03041223528485

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022431080653877
$ws.Range("D20").Value = 1.025351036197223
$ws.Range("E20").Value = 1.023235354378734
$ws.Range("F20").Value = 1.031786928174515
$ws.Range("I20").Value = 1.030839703446852
$ws.Range("J20").Value = 1.028719964171145
$ws.Range("K20").Value = 1.028762204850336
$ws.Range("L20").Value = 1.02665410169428
$ws.Range("M20").Value = 1.035175253233426
$ws.Range("N20").Value = 1.030180864516097

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021274600733309
$ws.Range("D21").Value = 1.024295790242747
$ws.Range("E21").Value = 1.022254434614573
$ws.Range("F21").Value = 1.030411116220063
$ws.Range("I21").Value = 1.030648212169285
$ws.Range("J21").Value = 1.027968219182023
$ws.Range("K21").Value = 1.027920038404071
$ws.Range("L21").Value = 1.025886490804826
$ws.Range("M21").Value = 1.034012178885277
$ws.Range("N21").Value = 1.029428051962864

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020547595062656
$ws.Range("D22").Value = 1.023632760318114
$ws.Range("E22").Value = 1.021638210841965
$ws.Range("F22").Value = 1.029546291860863
$ws.Range("I22").Value = 1.030526251626306
$ws.Range("J22").Value = 1.027495105701495
$ws.Range("K22").Value = 1.027390317722675
$ws.Range("L22").Value = 1.025403702408006
$ws.Range("M22").Value = 1.033280526522669
$ws.Range("N22").Value = 1.028954266606927

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020932983329952
$ws.Range("D23").Value = 1.02398420274432
$ws.Range("E23").Value = 1.021964833343976
$ws.Range("F23").Value = 1.030004732732312
$ws.Range("I23").Value = 1.030591053265496
$ws.Range("J23").Value = 1.027745955821637
$ws.Range("K23").Value = 1.027671153211795
$ws.Range("L23").Value = 1.025659652851516
$ws.Range("M23").Value = 1.033668425715196
$ws.Range("N23").Value = 1.029205472963012

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022450471236553
$ws.Range("D24").Value = 1.025368735061164
$ws.Range("E24").Value = 1.023251808402527
$ws.Range("F24").Value = 1.031809997374969
$ws.Range("I24").Value = 1.030842887327717
$ws.Range("J24").Value = 1.028732559443085
$ws.Range("K24").Value = 1.028776320188676
$ws.Range("L24").Value = 1.026666968047345
$ws.Range("M24").Value = 1.03519474592849
$ws.Range("N24").Value = 1.030193477674768

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024211960358649
$ws.Range("D25").Value = 1.026977336386326
$ws.Range("E25").Value = 1.024747518347874
$ws.Range("F25").Value = 1.033905829350285
$ws.Range("I25").Value = 1.031128318819796
$ws.Range("J25").Value = 1.029875444182363
$ws.Range("K25").Value = 1.030057853407662
$ws.Range("L25").Value = 1.027835192426218
$ws.Range("M25").Value = 1.036964323396776
$ws.Range("N25").Value = 1.031337985441467
